# Add a "Location" column (H) with "WareHouse" values for all data rows,
# mirroring the style used by the existing last column (G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 26

# Header cell (no special style, same as G1)
$ws.Range("H1").Value = "Location"

# Data cells - copy the format from column G (xlPasteFormats = -4122), then set values
$ws.Range("G2:G26").Copy()
$ws.Range("H2:H26").PasteSpecial(-4122)
$ws.Range("H2:H26").Value = "WareHouse"

# Column width for H
$ws.Columns.Item(8).ColumnWidth = 13.5703125

# Update selection to H1, matching the diff's sheetView change
$ws.Range("H1").Select()
